# Fruta / hortaliza, semanal
# Insert a new weekly record as row 24 in the data table, pushing the
# existing rows 24:81 down to 25:82 (dimension grows from A1:R81 to A1:R82).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 24:81 down by one (xlShiftDown = -4121) to make room for the
# new record at row 24. This also extends the used range / dimension.
$ws.Rows("24:24").Insert(-4121)

# Populate the newly inserted row 24 with the new weekly data point.
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 44868
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 100112026
$ws.Range("G24").Value = "Haba"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 90
$ws.Range("K24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = 10000
$ws.Range("N24").Value = "$/saco 25 kilos"
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 400
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
